$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "System, system, backup@backdoor.com" -> "System, backup@backdoor.com, system"
$rowsA = @(2, 28, 54)
foreach ($r in $rowsA) {
    $ws.Cells.Item($r, 7).Value = "System, backup@backdoor.com, system"
}

# Rows where "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
$rowsB = @(87, 113, 139)
foreach ($r in $rowsB) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, dnasr281@gmail.com"
}

# Rows where "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$rowsC = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)
foreach ($r in $rowsC) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

Write-Host "Applied all G-column Recorded By reorderings"
